# Add a new hike entry: "Owyhigh Lakes from Deer Creek"
# (10.5 miles, 3100 ft elevation gain, strenuous) into the sorted
# "Hike Difficulties" table, right before the alphabetically-later
# "Owyhigh Lakes from White River Road" entry, which currently sits
# at worksheet row 35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$insertRow = 35

# Shift row 35 (and everything below it) down by one row, opening up
# a blank row 35 for the new hike.
$ws.Rows.Item($insertRow).Insert()

# Populate the newly-opened row with the new hike's data.
$ws.Range("A" + $insertRow).Value = "Owyhigh Lakes from Deer Creek"
$ws.Range("B" + $insertRow).Value = 10.5
$ws.Range("C" + $insertRow).Value = 3100
$ws.Range("D" + $insertRow).Value = "strenuous"

# Grow Table1 so the structured table covers the new row too.
$lo = $ws.ListObjects.Item(1)
$newLastRow = $lo.Range.Row + $lo.Range.Rows.Count
$lo.Resize($ws.Range("A1:D" + $newLastRow))

# Match the author's final selection position.
$ws.Range("D36").Select()
